# Automatic update of files.
# Row 2 and Row 3 swap their "case" data (A, B, G, H, I, R and the
# S/T/V/W/X/Y HYPERLINK formulas), and every row's "Förändrad" date
# (column C, rows 2-9) moves from 46070 to 46072.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes what used to be Row 3 ("A 14042-2023") ---
$ws.Range("A2").Value = "A 14042-2023"
$ws.Range("B2").Value = 45008
$ws.Range("C2").Value = 46072
$ws.Range("G2").Value = 4.1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("R2").Value = "Revlummer`r`nÄkta lopplummer"

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/artfynd/A 14042-2023 artfynd.xlsx", "A 14042-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/kartor/A 14042-2023 karta.png", "A 14042-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomål/A 14042-2023 FSC-klagomål.docx", "A 14042-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomålsmail/A 14042-2023 FSC-klagomål mail.docx", "A 14042-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsyn/A 14042-2023 tillsynsbegäran.docx", "A 14042-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsynsmail/A 14042-2023 tillsynsbegäran mail.docx", "A 14042-2023")'

# --- Row 3 becomes what used to be Row 2 ("A 5398-2026") ---
$ws.Range("A3").Value = "A 5398-2026"
$ws.Range("B3").Value = 46050.49048611111
$ws.Range("C3").Value = 46072
$ws.Range("G3").Value = 0.5
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 2
$ws.Range("R3").Value = "Brandticka`r`nKambräken"

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/artfynd/A 5398-2026 artfynd.xlsx", "A 5398-2026")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/kartor/A 5398-2026 karta.png", "A 5398-2026")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomål/A 5398-2026 FSC-klagomål.docx", "A 5398-2026")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/klagomålsmail/A 5398-2026 FSC-klagomål mail.docx", "A 5398-2026")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsyn/A 5398-2026 tillsynsbegäran.docx", "A 5398-2026")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1402/tillsynsmail/A 5398-2026 tillsynsbegäran mail.docx", "A 5398-2026")'

# --- Rows 4-9: only "Förändrad" (col C) changes 46070 -> 46072 ---
$ws.Range("C4").Value = 46072
$ws.Range("C5").Value = 46072
$ws.Range("C6").Value = 46072
$ws.Range("C7").Value = 46072
$ws.Range("C8").Value = 46072
$ws.Range("C9").Value = 46072
